{"js": "const pairs = [\n  [\"92\u00d739=3588\", \"18\u00d718=324\"],\n  [\"98\u00d779=7742\", \"95\u00d788=8360\"],\n  [\"74\u00d760=4440\", \"14\u00d783=1162\"],\n  [\"83\u00d732=2656\", \"65\u00d796=6240\"],\n  [\"81\u00d753=4293\", \"87\u00d760=5220\"],\n  [\"20\u00d792=1840\", \"71\u00d756=3976\"],\n  [\"65\u00d754=3510\", \"15\u00d757=855\"],\n  [\"97\u00d735=3395\", \"76\u00d734=2584\"],\n  [\"60\u00d743=2580\", \"72\u00d763=4536\"],\n  [\"57\u00d753=3021\", \"60\u00d711=660\"],\n  [\"29\u00d715=435\", \"68\u00d774=5032\"],\n  [\"38\u00d746=1748\", \"63\u00d719=1197\"],\n  [\"91\u00d763=5733\", \"66\u00d729=1914\"],\n  [\"33\u00d761=2013\", \"49\u00d774=3626\"],\n  [\"38\u00d767=2546\", \"95\u00d756=5320\"],\n  [\"47\u00d788=4136\", \"42\u00d749=2058\"],\n  [\"65\u00d722=1430\", \"71\u00d722=1562\"],\n  [\"55\u00d735=1925\", \"91\u00d796=8736\"],\n  [\"56\u00d723=1288\", \"98\u00d778=7644\"],\n  [\"52\u00d760=3120\", \"26\u00d759=1534\"],\n  [\"43\u00d717=731\", \"78\u00d743=3354\"],\n  [\"33\u00d751=1683\", \"58\u00d770=4060\"],\n  [\"53\u00d731=1643\", \"51\u00d797=4947\"],\n  [\"77\u00d738=2926\", \"37\u00d734=1258\"],\n  [\"18\u00d768=1224\", \"37\u00d755=2035\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"92\u00d739=3588\", \"18\u00d718=324\"),\n    @(\"98\u00d779=7742\", \"95\u00d788=8360\"),\n    @(\"74\u00d760=4440\", \"14\u00d783=1162\"),\n    @(\"83\u00d732=2656\", \"65\u00d796=6240\"),\n    @(\"81\u00d753=4293\", \"87\u00d760=5220\"),\n    @(\"20\u00d792=1840\", \"71\u00d756=3976\"),\n    @(\"65\u00d754=3510\", \"15\u00d757=855\"),\n    @(\"97\u00d735=3395\", \"76\u00d734=2584\"),\n    @(\"60\u00d743=2580\", \"72\u00d763=4536\"),\n    @(\"57\u00d753=3021\", \"60\u00d711=660\"),\n    @(\"29\u00d715=435\", \"68\u00d774=5032\"),\n    @(\"38\u00d746=1748\", \"63\u00d719=1197\"),\n    @(\"91\u00d763=5733\", \"66\u00d729=1914\"),\n    @(\"33\u00d761=2013\", \"49\u00d774=3626\"),\n    @(\"38\u00d767=2546\", \"95\u00d756=5320\"),\n    @(\"47\u00d788=4136\", \"42\u00d749=2058\"),\n    @(\"65\u00d722=1430\", \"71\u00d722=1562\"),\n    @(\"55\u00d735=1925\", \"91\u00d796=8736\"),\n    @(\"56\u00d723=1288\", \"98\u00d778=7644\"),\n    @(\"52\u00d760=3120\", \"26\u00d759=1534\"),\n    @(\"43\u00d717=731\", \"78\u00d743=3354\"),\n    @(\"33\u00d751=1683\", \"58\u00d770=4060\"),\n    @(\"53\u00d731=1643\", \"51\u00d797=4947\"),\n    @(\"77\u00d738=2926\", \"37\u00d734=1258\"),\n    @(\"18\u00d768=1224\", \"37\u00d755=2035\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute([ref]$old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$new, 2)\n}\n"}
